$d = $word.ActiveDocument

# --- 1) Strike through the completed TODO items ---------------------------
# "Figure 1. Change "Distance" to "d""
$d.Paragraphs.Item(2).Range.Font.StrikeThrough = 1
# "Figure 2.  Add sketch of Polymer with loops corresponding to peaks"
$d.Paragraphs.Item(4).Range.Font.StrikeThrough = 1
# "Figure 2. Add encounter probability graph from simulations"
$d.Paragraphs.Item(5).Range.Font.StrikeThrough = 1

# --- 2) Update the Figure 3 formula text -----------------------------------
# old: Figure 3. Calculate the anomalous exponent <|x(t+dt)-x(t)|^2> for each bead in the TAD
# new: Figure 3. Calculate the anomalous exponent <|x(t)-x(0)|^2> for each bead in the TAD
$p7 = $d.Paragraphs.Item(7)
$p7start = $p7.Range.Start
$p7textRange = $d.Range($p7start, $p7.Range.End - 1)
$p7textRange.Text = "Figure 3. Calculate the anomalous exponent <|x(t)-x(0)|^2> for each bead in the TAD"

# Recreate the run split that existed around "t" / ")-x(0" by touching a
# collapsed range there (mirrors how Word breaks runs as you type/edit).
$splitPos = $p7start + 48
$rSplit = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TempSplit", $rSplit)
$d.Bookmarks.Item("TempSplit").Delete()

# Move the "_GoBack" bookmark (last-edit marker) from the end of the
# document to right after "...x(0" in this paragraph, i.e. between
# ")-x(0" and ")|^2...".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bmPos = $p7start + 53
$rBookmark = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $rBookmark)

# --- 3) Strike through the last TODO item ----------------------------------
# "Figure 4 add bars representing TAD D and E in the encounter histogram. "
$d.Paragraphs.Item(9).Range.Font.StrikeThrough = 1
